# Bug fix in Eduati data files: Sheet1 ("meas") had 43 stray extra rows
# (45-87) left over from a previous, longer dataset; trim it back down to
# the same 44 rows (1 header + 43 data rows) used by Sheet2/Sheet3, and
# leave the workbook's view state pointed at Sheet1 (where the edit was
# made) instead of Sheet3.

$wb = $excel.ActiveWorkbook

# --- 1. Trim the stray rows 45:87 off Sheet1 ------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A45:A87").EntireRow.Delete()

# --- 2. Leave the selection/active sheet where the fix was made ----------
# Sheet1 becomes the active/selected tab, scrolled down near the bottom of
# the trimmed data with E64 as the active cell.
$ws1.Activate()
$ws1.Range("E64").Select()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1

# Sheet3 (previously the active tab) goes back to its default top-left
# selection now that it is no longer the active sheet.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A2:N44").Select()

# Re-activate Sheet1 so it is left as the selected/active tab.
$ws1.Activate()
